$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 555.7143
$ws.Range("I18").Value = 548.3333
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 548.3333
$ws.Range("L18").Value = 600
$ws.Range("M18").Value = -264.3333
$ws.Range("N18").Value = -1168
$ws.Range("H64").Value = 4366
$ws.Range("I64").Value = 4285.636
$ws.Range("J64").Value = 4513.3335
$ws.Range("K64").Value = 4285.636
$ws.Range("L64").Value = 4513.3335
$ws.Range("M64").Value = -4037.636
$ws.Range("N64").Value = -5009.3335
$ws.Range("H67").Value = 4366
$ws.Range("I67").Value = 4285.636
$ws.Range("J67").Value = 4513.3335
$ws.Range("K67").Value = 4285.636
$ws.Range("L67").Value = 4513.3335
$ws.Range("M67").Value = -3427.636
$ws.Range("N67").Value = -6229.3335
$ws.Range("H76").Value = 3114.634
$ws.Range("I76").Value = 3072.7273
$ws.Range("K76").Value = 3072.7273
$ws.Range("M76").Value = -2757.7273
$ws.Range("H79").Value = 3114.634
$ws.Range("I79").Value = 3072.7273
$ws.Range("K79").Value = 3072.7273
$ws.Range("M79").Value = -1980.7273
$ws.Range("H114").Value = 44599
$ws.Range("J114").Value = 44599
$ws.Range("L114").Value = 44599
$ws.Range("N114").Value = -53277
$ws.Range("H129").Value = 859.3214
$ws.Range("J129").Value = 994.2174
$ws.Range("L129").Value = 2982.6522
$ws.Range("N129").Value = -12982.6522
$ws.Range("H131").Value = 4532.778
$ws.Range("I131").Value = 1798.3334
$ws.Range("J131").Value = 10001.667
$ws.Range("K131").Value = 5395.0002
$ws.Range("L131").Value = 30005.001
$ws.Range("M131").Value = -355.0002000000004
$ws.Range("N131").Value = -40085.001
$ws.Range("H132").Value = 3421.2415
$ws.Range("I132").Value = 3525.5715
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 10576.7145
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -8046.7145
$ws.Range("N132").Value = -6560
$ws.Range("H137").Value = 28573086
$ws.Range("I137").Value = 1020.16
$ws.Range("K137").Value = 3060.48
$ws.Range("M137").Value = -510.48

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1178686.1
$ws.Range("I2").Value = 1714
$ws.Range("J2").Value = 1472929.1
$ws.Range("K2").Value = 1714
$ws.Range("L2").Value = 1472929.1
$ws.Range("M2").Value = -1601
$ws.Range("N2").Value = -1473155.1
$ws.Range("H43").Value = 7325.2856
$ws.Range("J43").Value = 7325.2856
$ws.Range("L43").Value = 7325.2856
$ws.Range("N43").Value = -7951.2856
$ws.Range("H63").Value = 2473.625
$ws.Range("I63").Value = 1961.6364
$ws.Range("J63").Value = 3600
$ws.Range("K63").Value = 1961.6364
$ws.Range("L63").Value = 3600
$ws.Range("M63").Value = -1275.6364
$ws.Range("N63").Value = -4972
$ws.Range("H66").Value = 2473.625
$ws.Range("I66").Value = 1961.6364
$ws.Range("J66").Value = 3600
$ws.Range("K66").Value = 9808.182000000001
$ws.Range("L66").Value = 18000
$ws.Range("M66").Value = -6376.182000000001
$ws.Range("N66").Value = -24864
$ws.Range("H88").Value = 2584.2856
$ws.Range("I88").Value = 2850
$ws.Range("J88").Value = 2230
$ws.Range("K88").Value = 2850
$ws.Range("L88").Value = 2230
$ws.Range("M88").Value = -2444
$ws.Range("N88").Value = -3042
$ws.Range("H91").Value = 2584.2856
$ws.Range("I91").Value = 2850
$ws.Range("J91").Value = 2230
$ws.Range("K91").Value = 2850
$ws.Range("L91").Value = 2230
$ws.Range("M91").Value = -1446
$ws.Range("N91").Value = -5038
$ws.Range("H116").Value = 1178686.1
$ws.Range("I116").Value = 1714
$ws.Range("J116").Value = 1472929.1
$ws.Range("K116").Value = 1714
$ws.Range("L116").Value = 1472929.1
$ws.Range("M116").Value = 580
$ws.Range("N116").Value = -1477517.1
$ws.Range("H122").Value = 2313.0588
$ws.Range("I122").Value = 2196.5
$ws.Range("J122").Value = 2416.6667
$ws.Range("K122").Value = 6589.5
$ws.Range("L122").Value = 7250.000100000001
$ws.Range("M122").Value = -4139.5
$ws.Range("N122").Value = -12150.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1178686.1
$ws.Range("I3").Value = 1714
$ws.Range("J3").Value = 1472929.1
$ws.Range("K3").Value = 1714
$ws.Range("L3").Value = 1472929.1
$ws.Range("M3").Value = -1600
$ws.Range("N3").Value = -1473157.1
$ws.Range("H105").Value = 1754.0834
$ws.Range("I105").Value = 1727.6666
$ws.Range("K105").Value = 1727.6666
$ws.Range("M105").Value = 19.33339999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 11088
$ws.Range("J50").Value = 11088
$ws.Range("L50").Value = 11088
$ws.Range("N50").Value = -12338
$ws.Range("H51").Value = 9982.200000000001
$ws.Range("J51").Value = 10977.75
$ws.Range("L51").Value = 10977.75
$ws.Range("N51").Value = -12449.75
$ws.Range("H59").Value = 15732.714
$ws.Range("J59").Value = 15688.167
$ws.Range("L59").Value = 15688.167
$ws.Range("N59").Value = -17978.167
$ws.Range("H60").Value = 8831.200000000001
$ws.Range("J60").Value = 10164
$ws.Range("L60").Value = 10164
$ws.Range("N60").Value = -11186
$ws.Range("H61").Value = 9982.200000000001
$ws.Range("J61").Value = 10977.75
$ws.Range("L61").Value = 10977.75
$ws.Range("N61").Value = -11673.75
$ws.Range("H68").Value = 19066.5
$ws.Range("J68").Value = 19066.5
$ws.Range("L68").Value = 19066.5
$ws.Range("N68").Value = -20564.5
$ws.Range("H71").Value = 19066.5
$ws.Range("J71").Value = 19066.5
$ws.Range("L71").Value = 57199.5
$ws.Range("N71").Value = -64687.5
$ws.Range("H74").Value = 13801.4
$ws.Range("J74").Value = 16453.625
$ws.Range("L74").Value = 16453.625
$ws.Range("N74").Value = -18201.625
$ws.Range("H77").Value = 13801.4
$ws.Range("J77").Value = 16453.625
$ws.Range("L77").Value = 49360.875
$ws.Range("N77").Value = -58096.875
$ws.Range("H110").Value = 48902
$ws.Range("J110").Value = 48902
$ws.Range("L110").Value = 48902
$ws.Range("N110").Value = -57082

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 125305.94
$ws.Range("J107").Value = 111390.22
$ws.Range("L107").Value = 334170.66
$ws.Range("N107").Value = -338010.66
$ws.Range("H113").Value = 2708.9788
$ws.Range("I113").Value = 601.6667
$ws.Range("J113").Value = 4016.9656
$ws.Range("K113").Value = 1805.0001
$ws.Range("L113").Value = 12050.8968
$ws.Range("M113").Value = 364.9999
$ws.Range("N113").Value = -16390.8968
$ws.Range("H131").Value = 2045.8682
$ws.Range("I131").Value = 5929.909
$ws.Range("J131").Value = 1511.8125
$ws.Range("K131").Value = 17789.727
$ws.Range("L131").Value = 4535.4375
$ws.Range("M131").Value = -12749.727
$ws.Range("N131").Value = -14615.4375
$ws.Range("H139").Value = 1389.7106
$ws.Range("I139").Value = 896.5
$ws.Range("J139").Value = 2458.3333
$ws.Range("K139").Value = 2689.5
$ws.Range("L139").Value = 7374.999899999999
$ws.Range("M139").Value = 2450.5
$ws.Range("N139").Value = -17654.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7734.6895
$ws.Range("I70").Value = 10231.0625
$ws.Range("J70").Value = 4662.231
$ws.Range("K70").Value = 10231.0625
$ws.Range("L70").Value = 4662.231
$ws.Range("M70").Value = -9961.0625
$ws.Range("N70").Value = -5202.231
$ws.Range("H73").Value = 7734.6895
$ws.Range("I73").Value = 10231.0625
$ws.Range("J73").Value = 4662.231
$ws.Range("K73").Value = 10231.0625
$ws.Range("L73").Value = 4662.231
$ws.Range("M73").Value = -9295.0625
$ws.Range("N73").Value = -6534.231
$ws.Range("H80").Value = 28250
$ws.Range("I80").Value = 7000
$ws.Range("J80").Value = 35333.332
$ws.Range("K80").Value = 7000
$ws.Range("L80").Value = 35333.332
$ws.Range("M80").Value = -6002
$ws.Range("N80").Value = -37329.332
$ws.Range("H83").Value = 28250
$ws.Range("I83").Value = 7000
$ws.Range("J83").Value = 35333.332
$ws.Range("K83").Value = 35000
$ws.Range("L83").Value = 176666.66
$ws.Range("M83").Value = -30008
$ws.Range("N83").Value = -186650.66

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1756.7778
$ws.Range("I122").Value = 1701.0625
$ws.Range("J122").Value = 2202.5
$ws.Range("K122").Value = 5103.1875
$ws.Range("L122").Value = 6607.5
$ws.Range("M122").Value = -2653.1875
$ws.Range("N122").Value = -11507.5
